$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config")

# Update the path values to support multiple tables with the same mail
$ws.Range("B7").Value = "Data/HTML.txt"
$ws.Range("B8").Value = "Output/"

# Update the active selection to E9
$ws.Range("E9").Select()
